$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Snippets")

# The "fields" snippet id was renamed from "word-manage-fields" to
# "word-document-manage-fields" (design change for GA promotion).
# This value lives in column E (SnippetIdIntheYAMLFile) for the rows
# describing Body.fields / Field.* / FieldCollection.* members.
$oldValue = "word-manage-fields"
$newValue = "word-document-manage-fields"

$rowsToUpdate = @(9, 37, 38, 39, 40, 41)
foreach ($r in $rowsToUpdate) {
    $cell = $ws.Cells.Item($r, 5)
    if ($cell.Value() -eq $oldValue) {
        $cell.Value = $newValue
    }
}

# Reflect the author's last cell selection recorded in the saved file.
$ws.Range("E10").Select()
